{"js": "// \"L\u00e0m tr\u00f2n s\u1ed1 ti\u1ec1n\" edit: the \"B\u00ecnh qu\u00e2n v\u00f9ng\" (regional average) clause in\n// four paragraphs reads \"<field> % (<field>)\" even though that clause is not\n// actually a percentage there \u2014 three of them report money amounts (\"\u0111\u1ed3ng\")\n// and one reports a day count (\"ng\u00e0y\"). The fix rewrites the literal\n// \" % (\" text that sits between the two field placeholders to \" ng\u00e0y (\" /\n// \" \u0111\u1ed3ng (\" as appropriate, leaving the red field-code runs (the {X..}\n// placeholders on either side) completely untouched.\n//\n// Each target paragraph is identified by a field placeholder that is unique\n// in the document; within that single paragraph, the literal \" % (\" text is\n// also unique, so searching scoped to the paragraph (rather than the whole\n// body) unambiguously selects just the run/text we need to change without\n// touching the neighboring field-code runs' formatting.\n\nconst targets = [\n  { anchorField: \"{X16}\", word: \"ng\u00e0y\" },\n  { anchorField: \"{X23}\", word: \"\u0111\u1ed3ng\" },\n  { anchorField: \"{X30}\", word: \"\u0111\u1ed3ng\" },\n  { anchorField: \"{X37}\", word: \"\u0111\u1ed3ng\" }\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const { anchorField, word } of targets) {\n  const paragraph = paragraphs.items.find((p) => p.text.includes(anchorField));\n  if (!paragraph) {\n    throw new Error(\"Could not find paragraph containing \" + anchorField);\n  }\n\n  const hits = paragraph.search(\" % (\", { matchCase: true });\n  hits.load(\"text\");\n  await context.sync();\n\n  if (hits.items.length === 0) {\n    throw new Error('Could not find \" % (\" near ' + anchorField);\n  }\n\n  hits.items[0].insertText(\" \" + word + \" (\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"L\u00e0m tr\u00f2n s\u1ed1 ti\u1ec1n\" edit: the \"B\u00ecnh qu\u00e2n v\u00f9ng\" (regional average) clause in\n# four paragraphs reads \"<field> % (<field>)\" even though that clause is not\n# actually a percentage there \u2014 three of them report money amounts (\"\u0111\u1ed3ng\")\n# and one reports a day count (\"ng\u00e0y\"). The fix rewrites the literal\n# \" % (\" text that sits between the two field placeholders to \" ng\u00e0y (\" /\n# \" \u0111\u1ed3ng (\" as appropriate, leaving the red field-code runs (the {X..}\n# placeholders on either side) completely untouched.\n#\n# Each target spot is located by first finding the (document-unique) field\n# placeholder that precedes it, collapsing the range to right after that\n# placeholder, and then finding \" % (\" from there \u2014 this keeps the edit\n# scoped to just those four characters instead of touching the neighboring\n# field-code runs' formatting.\n\n$d = $word.ActiveDocument\n\n$targets = @(\n  @{ Anchor = \"{X16}\"; Word = \"ng\u00e0y\" },\n  @{ Anchor = \"{X23}\"; Word = \"\u0111\u1ed3ng\" },\n  @{ Anchor = \"{X30}\"; Word = \"\u0111\u1ed3ng\" },\n  @{ Anchor = \"{X37}\"; Word = \"\u0111\u1ed3ng\" }\n)\n\nforeach ($t in $targets) {\n  $rng = $d.Content\n\n  $foundAnchor = $rng.Find.Execute($t.Anchor)\n  if (-not $foundAnchor) {\n    throw \"Could not find anchor text: $($t.Anchor)\"\n  }\n\n  $rng.Collapse(0)  # wdCollapseEnd: collapse to right after the field placeholder\n\n  $foundSuffix = $rng.Find.Execute(\" % (\")\n  if (-not $foundSuffix) {\n    throw \"Could not find `\" % (`\" after $($t.Anchor)\"\n  }\n\n  $rng.Text = \" \" + $t.Word + \" (\"\n}\n"}
